$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing the existing rows 11-13 down to rows 12-14.
$ws.Rows("11:11").Insert()

# Populate the new row 11 with the new asserted_distribution test record:
#   otu_name = albert, state = Nebraska, citation = "philbert in Nebraska as described by Anon."
$ws.Range("A11").Value = "albert"
$ws.Range("E11").Value = "Nebraska"
$ws.Range("I11").Value = "philbert in Nebraska as described by Anon."

# The new row uses the taller row height seen on similar citation rows.
$ws.Rows("11:11").RowHeight = 75

# Reflect the new selection/scroll position left behind by the edit.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I11").Select()
